# Update "想去人数" (F column) counts per the diff for each sheet
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(2, 6).Value = 1731
$ws.Cells.Item(3, 6).Value = 10105
$ws.Cells.Item(5, 6).Value = 15
$ws.Cells.Item(6, 6).Value = 593
$ws.Cells.Item(7, 6).Value = 72
$ws.Cells.Item(8, 6).Value = 1600
$ws.Cells.Item(9, 6).Value = 168
$ws.Cells.Item(10, 6).Value = 373
$ws.Cells.Item(15, 6).Value = 1172
$ws.Cells.Item(16, 6).Value = 127
$ws.Cells.Item(18, 6).Value = 11
$ws.Cells.Item(19, 6).Value = 86
$ws.Cells.Item(20, 6).Value = 344
$ws.Cells.Item(21, 6).Value = 14
$ws.Cells.Item(22, 6).Value = 316
$ws.Cells.Item(23, 6).Value = 100
$ws.Cells.Item(24, 6).Value = 1154
$ws.Cells.Item(25, 6).Value = 694
$ws.Cells.Item(26, 6).Value = 21
$ws.Cells.Item(27, 6).Value = 37
$ws.Cells.Item(29, 6).Value = 224
$ws.Cells.Item(31, 6).Value = 392
$ws.Cells.Item(32, 6).Value = 216
$ws.Cells.Item(33, 6).Value = 370
$ws.Cells.Item(34, 6).Value = 525
$ws.Cells.Item(35, 6).Value = 599
$ws.Cells.Item(36, 6).Value = 734
$ws.Cells.Item(37, 6).Value = 525
$ws.Cells.Item(38, 6).Value = 1268
$ws.Cells.Item(39, 6).Value = 810
$ws.Cells.Item(40, 6).Value = 381
$ws.Cells.Item(42, 6).Value = 7
$ws.Cells.Item(44, 6).Value = 75
$ws.Cells.Item(45, 6).Value = 352

$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(4, 6).Value = 45
$ws.Cells.Item(9, 6).Value = 4
$ws.Cells.Item(11, 6).Value = 56
$ws.Cells.Item(14, 6).Value = 98
$ws.Cells.Item(15, 6).Value = 59
$ws.Cells.Item(18, 6).Value = 1083
$ws.Cells.Item(20, 6).Value = 592
$ws.Cells.Item(21, 6).Value = 1099
$ws.Cells.Item(24, 6).Value = 74
$ws.Cells.Item(29, 6).Value = 9
$ws.Cells.Item(31, 6).Value = 205
$ws.Cells.Item(35, 6).Value = 185
$ws.Cells.Item(36, 6).Value = 43
$ws.Cells.Item(41, 6).Value = 64
$ws.Cells.Item(42, 6).Value = 37

$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(4, 6).Value = 802
$ws.Cells.Item(5, 6).Value = 190
$ws.Cells.Item(6, 6).Value = 2511
$ws.Cells.Item(7, 6).Value = 4052
$ws.Cells.Item(8, 6).Value = 57
$ws.Cells.Item(10, 6).Value = 285
$ws.Cells.Item(11, 6).Value = 185

$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(2, 6).Value = 1731
$ws.Cells.Item(3, 6).Value = 802
$ws.Cells.Item(4, 6).Value = 10105
$ws.Cells.Item(5, 6).Value = 190
$ws.Cells.Item(7, 6).Value = 4052
$ws.Cells.Item(8, 6).Value = 57
$ws.Cells.Item(9, 6).Value = 285
$ws.Cells.Item(10, 6).Value = 286
$ws.Cells.Item(11, 6).Value = 593
$ws.Cells.Item(12, 6).Value = 1600
$ws.Cells.Item(13, 6).Value = 168
$ws.Cells.Item(14, 6).Value = 373
$ws.Cells.Item(16, 6).Value = 4
$ws.Cells.Item(18, 6).Value = 1172
$ws.Cells.Item(19, 6).Value = 127
$ws.Cells.Item(20, 6).Value = 56
$ws.Cells.Item(22, 6).Value = 98
$ws.Cells.Item(23, 6).Value = 59
$ws.Cells.Item(24, 6).Value = 86
$ws.Cells.Item(25, 6).Value = 1083
$ws.Cells.Item(26, 6).Value = 344
$ws.Cells.Item(27, 6).Value = 316
$ws.Cells.Item(28, 6).Value = 1099
$ws.Cells.Item(30, 6).Value = 1154
$ws.Cells.Item(31, 6).Value = 694
$ws.Cells.Item(32, 6).Value = 74
$ws.Cells.Item(35, 6).Value = 392
$ws.Cells.Item(36, 6).Value = 9
$ws.Cells.Item(37, 6).Value = 370
$ws.Cells.Item(38, 6).Value = 525
$ws.Cells.Item(39, 6).Value = 599
$ws.Cells.Item(40, 6).Value = 205
$ws.Cells.Item(41, 6).Value = 734
$ws.Cells.Item(42, 6).Value = 525
$ws.Cells.Item(43, 6).Value = 810
$ws.Cells.Item(44, 6).Value = 381
$ws.Cells.Item(45, 6).Value = 43
$ws.Cells.Item(49, 6).Value = 352
$ws.Cells.Item(50, 6).Value = 64
